$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "nan" value that was in C6 (row 6 keeps its other original data)
$ws.Range("C6").Value = ""

# Insert a new row 7 with the new review data
$ws.Range("A7").Value = "parisk"
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = "nan"
$ws.Range("D7").Value = "DIS"
$ws.Range("E7").Value = "WRI"
$ws.Range("F7").Value = "dc9804e9-fe90-49ab-88bb-ac97478c1b97"
$ws.Range("G7").Value = "i87JIQTAnB8AQ_annotated.xlsx"
$ws.Range("H7").Value = "As you suggested, I did run comparison tests and I will present the results here."
$ws.Range("I7").Value = "Correct"
